# Auto update: 2025-05-20 19:20:14
# Reorders the company entries (columns A, B, H) across rows 2-8 while
# leaving the other per-row metadata (C, D, E, F, G) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns A (Company Name), B (Company Number), H (Category)
# for each affected row, reflecting the reordering described in the diff.
$rows = @{
    2 = @("T GILPIN PHYSIO CONSULTANCY LTD", "16460503", "LP")
    3 = @("SAMVIV PARTNERS LTD", "16460672", "Partners")
    4 = @("4D CAPITAL PROPCO (44) LIMITED", "16461269", "Capital")
    6 = @("AFROSCOT VENTURES LTD", "16462878", "Ventures")
    7 = @("ST GEORGE CAPITAL (LAND) LIMITED", "16462880", "Capital")
    8 = @("DGPI LTD", "SC849118", "GP")
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    # Company Number must stay text (the source keeps these as plain
    # strings even when they look numeric), so force text format before
    # writing the value to avoid Excel auto-converting it to a number.
    $ws.Range("B$r").NumberFormat = "@"
    $ws.Range("A$r").Value = $vals[0]
    $ws.Range("B$r").Value = $vals[1]
    $ws.Range("H$r").Value = $vals[2]
}
